$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3825.524
$ws.Range("I76").Value = 3560
$ws.Range("J76").Value = 3908.5
$ws.Range("K76").Value = 3560
$ws.Range("L76").Value = 3908.5
$ws.Range("M76").Value = -3245
$ws.Range("N76").Value = -4538.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3825.524
$ws.Range("I79").Value = 3560
$ws.Range("J79").Value = 3908.5
$ws.Range("K79").Value = 3560
$ws.Range("L79").Value = 3908.5
$ws.Range("M79").Value = -2468
$ws.Range("N79").Value = -6092.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1853397.5
$ws.Range("I137").Value = 2273722
$ws.Range("J137").Value = 3969.9
$ws.Range("K137").Value = 6821166
$ws.Range("L137").Value = 11909.7
$ws.Range("M137").Value = -6818616
$ws.Range("N137").Value = -17009.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1147.5
$ws.Range("J4").Value = 1263.3334
$ws.Range("L4").Value = 1263.3334
$ws.Range("N4").Value = -1495.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22191.977
$ws.Range("I32").Value = 21142.393
$ws.Range("J32").Value = 25102.182
$ws.Range("K32").Value = 21142.393
$ws.Range("L32").Value = 25102.182
$ws.Range("M32").Value = -20855.393
$ws.Range("N32").Value = -25676.182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1220.3077
$ws.Range("I45").Value = 1047.4
$ws.Range("J45").Value = 1796.6666
$ws.Range("K45").Value = 1047.4
$ws.Range("L45").Value = 1796.6666
$ws.Range("M45").Value = -670.4000000000001
$ws.Range("N45").Value = -2550.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20041248
$ws.Range("I61").Value = 23834350
$ws.Range("K61").Value = 23834350
$ws.Range("M61").Value = -23834138

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7622432.5
$ws.Range("I74").Value = 9834368
$ws.Range("K74").Value = 9834368
$ws.Range("M74").Value = -9833494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7622432.5
$ws.Range("I77").Value = 9834368
$ws.Range("K77").Value = 49171840
$ws.Range("M77").Value = -49167472

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45461.81
$ws.Range("I132").Value = 34062.766
$ws.Range("J132").Value = 65577.766
$ws.Range("K132").Value = 102188.298
$ws.Range("L132").Value = 196733.298
$ws.Range("M132").Value = -99658.29800000001
$ws.Range("N132").Value = -201793.298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 20041248
$ws.Range("I136").Value = 23834350
$ws.Range("K136").Value = 71503050
$ws.Range("M136").Value = -71500500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 31996.666
$ws.Range("J62").Value = 31996.666
$ws.Range("L62").Value = 31996.666
$ws.Range("N62").Value = -33368.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 31996.666
$ws.Range("J65").Value = 31996.666
$ws.Range("L65").Value = 95989.99800000001
$ws.Range("N65").Value = -102853.998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 25833.334
$ws.Range("J109").Value = 25833.334
$ws.Range("L109").Value = 25833.334
$ws.Range("N109").Value = -28607.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1888.4419
$ws.Range("I31").Value = 1063.5834
$ws.Range("J31").Value = 6130.5713
$ws.Range("K31").Value = 1063.5834
$ws.Range("L31").Value = 6130.5713
$ws.Range("M31").Value = -768.5834
$ws.Range("N31").Value = -6720.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1888.4419
$ws.Range("I34").Value = 1063.5834
$ws.Range("J34").Value = 6130.5713
$ws.Range("K34").Value = 1063.5834
$ws.Range("L34").Value = 6130.5713
$ws.Range("M34").Value = -861.5834
$ws.Range("N34").Value = -6534.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17545370
$ws.Range("I58").Value = 23257286
$ws.Range("K58").Value = 23257286
$ws.Range("M58").Value = -23257083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 30733.742
$ws.Range("I132").Value = 2204.9048
$ws.Range("K132").Value = 6614.714399999999
$ws.Range("M132").Value = -4084.714399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 27274.88
$ws.Range("I134").Value = 1531.697
$ws.Range("J134").Value = 121666.555
$ws.Range("K134").Value = 4595.090999999999
$ws.Range("L134").Value = 364999.665
$ws.Range("M134").Value = -2060.090999999999
$ws.Range("N134").Value = -370069.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 17545370
$ws.Range("I136").Value = 23257286
$ws.Range("K136").Value = 69771858
$ws.Range("M136").Value = -69769308

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 59000
$ws.Range("J140").Value = 59000
$ws.Range("L140").Value = 59000
$ws.Range("N140").Value = -69360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 4500
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -4846

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2842
$ws.Range("I63").Value = 2302.5
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 6907.5
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -6158.5
$ws.Range("N63").Value = -16498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3500.1
$ws.Range("I64").Value = 1612
$ws.Range("K64").Value = 4836
$ws.Range("M64").Value = -4566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 2842
$ws.Range("I66").Value = 2302.5
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 20722.5
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -16978.5
$ws.Range("N66").Value = -52488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3500.1
$ws.Range("I67").Value = 1612
$ws.Range("K67").Value = 4836
$ws.Range("M67").Value = -3900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 778.6667
$ws.Range("J98").Value = 894.4
$ws.Range("L98").Value = 2683.2
$ws.Range("N98").Value = -5679.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 48529704
$ws.Range("I121").Value = 1832.5
$ws.Range("J121").Value = 54238868
$ws.Range("K121").Value = 5497.5
$ws.Range("L121").Value = 162716604
$ws.Range("M121").Value = -4187.5
$ws.Range("N121").Value = -162719224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1799
$ws.Range("I132").Value = 788.5714
$ws.Range("J132").Value = 2683.125
$ws.Range("K132").Value = 7097.1426
$ws.Range("L132").Value = 24148.125
$ws.Range("M132").Value = -4567.1426
$ws.Range("N132").Value = -29208.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3306.25
$ws.Range("I138").Value = 2610
$ws.Range("J138").Value = 4466.6665
$ws.Range("K138").Value = 7830
$ws.Range("L138").Value = 13399.9995
$ws.Range("M138").Value = -2690
$ws.Range("N138").Value = -23679.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3669.4092
$ws.Range("I80").Value = 3163.125
$ws.Range("J80").Value = 3958.7144
$ws.Range("K80").Value = 3163.125
$ws.Range("L80").Value = 3958.7144
$ws.Range("M80").Value = -2165.125
$ws.Range("N80").Value = -5954.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3669.4092
$ws.Range("I83").Value = 3163.125
$ws.Range("J83").Value = 3958.7144
$ws.Range("K83").Value = 15815.625
$ws.Range("L83").Value = 19793.572
$ws.Range("M83").Value = -10823.625
$ws.Range("N83").Value = -29777.572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1800
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1800
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2394.2
$ws.Range("I126").Value = 1550
$ws.Range("J126").Value = 2957
$ws.Range("K126").Value = 4650
$ws.Range("L126").Value = 8871
$ws.Range("M126").Value = -2180
$ws.Range("N126").Value = -13811

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5749.6665
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 7624.5
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 7624.5
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -7848.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3047.75
$ws.Range("I40").Value = 2563.4583
$ws.Range("J40").Value = 4500.625
$ws.Range("K40").Value = 2563.4583
$ws.Range("L40").Value = 4500.625
$ws.Range("M40").Value = -2427.4583
$ws.Range("N40").Value = -4772.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1789.6316
$ws.Range("I100").Value = 1557.5714
$ws.Range("K100").Value = 1557.5714
$ws.Range("M100").Value = -1016.5714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5749.6665
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 7624.5
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 22873.5
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -27813.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 21114.098
$ws.Range("I132").Value = 1071
$ws.Range("J132").Value = 79701.62
$ws.Range("K132").Value = 3213
$ws.Range("L132").Value = 239104.86
$ws.Range("M132").Value = -683
$ws.Range("N132").Value = -244164.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 45540.934
$ws.Range("I136").Value = 25843.55
$ws.Range("J136").Value = 203120
$ws.Range("K136").Value = 77530.64999999999
$ws.Range("L136").Value = 609360
$ws.Range("M136").Value = -74980.64999999999
$ws.Range("N136").Value = -614460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 63659.25
$ws.Range("I100").Value = 167117.33
$ws.Range("J100").Value = 39784.31
$ws.Range("K100").Value = 334234.66
$ws.Range("L100").Value = 79568.62
$ws.Range("M100").Value = -333693.66
$ws.Range("N100").Value = -80650.62

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1975.2941
$ws.Range("I122").Value = 1295.0333
$ws.Range("J122").Value = 2947.0952
$ws.Range("K122").Value = 3885.0999
$ws.Range("L122").Value = 8841.285600000001
$ws.Range("M122").Value = -1435.0999
$ws.Range("N122").Value = -13741.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35526.516
$ws.Range("I132").Value = 23969.883
$ws.Range("J132").Value = 68655.53
$ws.Range("K132").Value = 71909.649
$ws.Range("L132").Value = 205966.59
$ws.Range("M132").Value = -69379.649
$ws.Range("N132").Value = -211026.59
